# Update annotations for Ying Tang
#
# The "cr dataset" annotation that used to live in row 36 is kept in row 36
# (its politeness_score in column B becomes a genuine number instead of
# text), and a brand-new annotation row (about competing/augmenting
# subgraph effects) is appended as row 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: politeness_score becomes numeric 3 (was text "3")
$ws.Cells.Item(36, 2).Value = 3

# Row 37: new annotation entry
$ws.Cells.Item(37, 1).Value = "Ying Tang"

# politeness_score stored as text "3" (matches the rest of the sheet's
# convention of inline-string "3" rather than a number). Force text type
# via NumberFormat, then strip the formatting residue so the cell keeps
# the workbook's default (unstyled) look.
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = "3"
$ws.Cells.Item(37, 2).ClearFormats()

$ws.Cells.Item(37, 3).Value = "无"
$ws.Cells.Item(37, 4).Value = "CRT"
$ws.Cells.Item(37, 5).Value = "MET"
$ws.Cells.Item(37, 6).Value = "b3917550-3902-443d-ae6f-4c206bcc883a"
$ws.Cells.Item(37, 7).Value = "HkJ1rgbCb_annotated.xlsx"
$ws.Cells.Item(37, 8).Value = "However, these selections do not seem to directly incorporate the competing/augmenting effects of having different subgraphs within a molecule."
